$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")

# Switch the selected tile material from Stone to Glass (drives the whole
# recalculation cascade: Variables sheet strings, HLOOKUPs, and charts 1 & 2).
$ws.Range("H2").Value = "Glass"

# Look up AtkPnt "per swing" instead of "per second" for the L column.
$ws.Range("L2").Value = "AtkPnt per swing"

# Rename the "Can be buffed" label to "Actorvalue".
$ws.Range("C3").Value = "Actorvalue"

# Add a new column header "Wut" and its value -1.
$ws.Range("J3").Value = "Wut"
$ws.Range("J4").Value = -1

# Update the selection to match the author's final cursor position.
$ws.Range("A4:C4").Select()

$wb.Save()
